# Add the new "2022-Q4" sheet right after "总计" (the first sheet),
# matching the position the other quarterly sheets get pushed down to.
$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $zj)
$q4.Name = "2022-Q4"

# ---- Header row (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名) ----
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Copy header styling (bold / centered / bordered) from the "总计" sheet's own
# header cell so the look matches the rest of the workbook.
$zj.Range("B1").Copy() | Out-Null
$q4.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$q4.Range("B1").Value = "基金代码"

# ---- Data rows: columns B:G hold text (fund codes keep leading zeros,
# percentages/amounts keep their printed decimal places), column A and H
# are numeric (row index / rank). ----
$q4.Range("B2:G29").NumberFormat = "@"

$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "512660"
$q4.Cells.Item(2,3).Value = "国泰中证军工ETF"
$q4.Cells.Item(2,4).Value = "104.53"
$q4.Cells.Item(2,5).Value = "99.72"
$q4.Cells.Item(2,6).Value = "2.53"
$q4.Cells.Item(2,7).Value = "2.6446"
$q4.Cells.Item(2,8).Value = 10
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "161024"
$q4.Cells.Item(3,3).Value = "富国中证军工指数A"
$q4.Cells.Item(3,4).Value = "50.37"
$q4.Cells.Item(3,5).Value = "94.33"
$q4.Cells.Item(3,6).Value = "2.32"
$q4.Cells.Item(3,7).Value = "1.1686"
$q4.Cells.Item(3,8).Value = 10
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "512680"
$q4.Cells.Item(4,3).Value = "广发中证军工ETF"
$q4.Cells.Item(4,4).Value = "32.36"
$q4.Cells.Item(4,5).Value = "99.61"
$q4.Cells.Item(4,6).Value = "2.46"
$q4.Cells.Item(4,7).Value = "0.7961"
$q4.Cells.Item(4,8).Value = 10
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "001556"
$q4.Cells.Item(5,3).Value = "天弘中证500指数增强A"
$q4.Cells.Item(5,4).Value = "25.50"
$q4.Cells.Item(5,5).Value = "94.27"
$q4.Cells.Item(5,6).Value = "1.91"
$q4.Cells.Item(5,7).Value = "0.4870"
$q4.Cells.Item(5,8).Value = 1
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "000478"
$q4.Cells.Item(6,3).Value = "建信中证500指数增强A"
$q4.Cells.Item(6,4).Value = "48.60"
$q4.Cells.Item(6,5).Value = "83.82"
$q4.Cells.Item(6,6).Value = "0.93"
$q4.Cells.Item(6,7).Value = "0.4520"
$q4.Cells.Item(6,8).Value = 10
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "001557"
$q4.Cells.Item(7,3).Value = "天弘中证500指数增强C"
$q4.Cells.Item(7,4).Value = "13.20"
$q4.Cells.Item(7,5).Value = "94.27"
$q4.Cells.Item(7,6).Value = "1.91"
$q4.Cells.Item(7,7).Value = "0.2521"
$q4.Cells.Item(7,8).Value = 1
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = "163115"
$q4.Cells.Item(8,3).Value = "申万菱信中证军工指数A"
$q4.Cells.Item(8,4).Value = "8.35"
$q4.Cells.Item(8,5).Value = "93.84"
$q4.Cells.Item(8,6).Value = "2.29"
$q4.Cells.Item(8,7).Value = "0.1912"
$q4.Cells.Item(8,8).Value = 10
$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).Value = "512560"
$q4.Cells.Item(9,3).Value = "易方达中证军工ETF"
$q4.Cells.Item(9,4).Value = "6.50"
$q4.Cells.Item(9,5).Value = "99.17"
$q4.Cells.Item(9,6).Value = "2.44"
$q4.Cells.Item(9,7).Value = "0.1586"
$q4.Cells.Item(9,8).Value = 10
$q4.Cells.Item(10,1).Value = 8
$q4.Cells.Item(10,2).Value = "502003"
$q4.Cells.Item(10,3).Value = "易方达军工指数（LOF）A"
$q4.Cells.Item(10,4).Value = "6.64"
$q4.Cells.Item(10,5).Value = "94.73"
$q4.Cells.Item(10,6).Value = "2.33"
$q4.Cells.Item(10,7).Value = "0.1547"
$q4.Cells.Item(10,8).Value = 10
$q4.Cells.Item(11,1).Value = 9
$q4.Cells.Item(11,2).Value = "003857"
$q4.Cells.Item(11,3).Value = "前海开源周期优选灵活配置混合A"
$q4.Cells.Item(11,4).Value = "2.13"
$q4.Cells.Item(11,5).Value = "89.59"
$q4.Cells.Item(11,6).Value = "5.69"
$q4.Cells.Item(11,7).Value = "0.1212"
$q4.Cells.Item(11,8).Value = 3
$q4.Cells.Item(12,1).Value = 10
$q4.Cells.Item(12,2).Value = "512810"
$q4.Cells.Item(12,3).Value = "华宝中证军工ETF"
$q4.Cells.Item(12,4).Value = "4.63"
$q4.Cells.Item(12,5).Value = "98.53"
$q4.Cells.Item(12,6).Value = "2.41"
$q4.Cells.Item(12,7).Value = "0.1116"
$q4.Cells.Item(12,8).Value = 10
$q4.Cells.Item(13,1).Value = 11
$q4.Cells.Item(13,2).Value = "159610"
$q4.Cells.Item(13,3).Value = "景顺长城中证500增强策略ETF"
$q4.Cells.Item(13,4).Value = "6.09"
$q4.Cells.Item(13,5).Value = "98.72"
$q4.Cells.Item(13,6).Value = "1.43"
$q4.Cells.Item(13,7).Value = "0.0871"
$q4.Cells.Item(13,8).Value = 1
$q4.Cells.Item(14,1).Value = 12
$q4.Cells.Item(14,2).Value = "002076"
$q4.Cells.Item(14,3).Value = "浙商中证500指数增强A"
$q4.Cells.Item(14,4).Value = "6.56"
$q4.Cells.Item(14,5).Value = "87.04"
$q4.Cells.Item(14,6).Value = "1.23"
$q4.Cells.Item(14,7).Value = "0.0807"
$q4.Cells.Item(14,8).Value = 1
$q4.Cells.Item(15,1).Value = 13
$q4.Cells.Item(15,2).Value = "003858"
$q4.Cells.Item(15,3).Value = "前海开源周期优选灵活配置混合C"
$q4.Cells.Item(15,4).Value = "0.72"
$q4.Cells.Item(15,5).Value = "89.59"
$q4.Cells.Item(15,6).Value = "5.69"
$q4.Cells.Item(15,7).Value = "0.0410"
$q4.Cells.Item(15,8).Value = 3
$q4.Cells.Item(16,1).Value = 14
$q4.Cells.Item(16,2).Value = "012842"
$q4.Cells.Item(16,3).Value = "易方达军工指数（LOF）C"
$q4.Cells.Item(16,4).Value = "1.36"
$q4.Cells.Item(16,5).Value = "94.73"
$q4.Cells.Item(16,6).Value = "2.33"
$q4.Cells.Item(16,7).Value = "0.0317"
$q4.Cells.Item(16,8).Value = 10
$q4.Cells.Item(17,1).Value = 15
$q4.Cells.Item(17,2).Value = "005633"
$q4.Cells.Item(17,3).Value = "建信中证500指数增强C"
$q4.Cells.Item(17,4).Value = "3.38"
$q4.Cells.Item(17,5).Value = "83.82"
$q4.Cells.Item(17,6).Value = "0.93"
$q4.Cells.Item(17,7).Value = "0.0314"
$q4.Cells.Item(17,8).Value = 10
$q4.Cells.Item(18,1).Value = 16
$q4.Cells.Item(18,2).Value = "013035"
$q4.Cells.Item(18,3).Value = "富国中证军工指数C"
$q4.Cells.Item(18,4).Value = "1.33"
$q4.Cells.Item(18,5).Value = "94.33"
$q4.Cells.Item(18,6).Value = "2.32"
$q4.Cells.Item(18,7).Value = "0.0309"
$q4.Cells.Item(18,8).Value = 10
$q4.Cells.Item(19,1).Value = 17
$q4.Cells.Item(19,2).Value = "002316"
$q4.Cells.Item(19,3).Value = "创金合信中证500指数增强C"
$q4.Cells.Item(19,4).Value = "2.72"
$q4.Cells.Item(19,5).Value = "93.95"
$q4.Cells.Item(19,6).Value = "1.10"
$q4.Cells.Item(19,7).Value = "0.0299"
$q4.Cells.Item(19,8).Value = 6
$q4.Cells.Item(20,1).Value = 18
$q4.Cells.Item(20,2).Value = "002311"
$q4.Cells.Item(20,3).Value = "创金合信中证500指数增强A"
$q4.Cells.Item(20,4).Value = "2.66"
$q4.Cells.Item(20,5).Value = "93.95"
$q4.Cells.Item(20,6).Value = "1.10"
$q4.Cells.Item(20,7).Value = "0.0293"
$q4.Cells.Item(20,8).Value = 6
$q4.Cells.Item(21,1).Value = 19
$q4.Cells.Item(21,2).Value = "007386"
$q4.Cells.Item(21,3).Value = "浙商中证500指数增强C"
$q4.Cells.Item(21,4).Value = "1.70"
$q4.Cells.Item(21,5).Value = "87.04"
$q4.Cells.Item(21,6).Value = "1.23"
$q4.Cells.Item(21,7).Value = "0.0209"
$q4.Cells.Item(21,8).Value = 1
$q4.Cells.Item(22,1).Value = 20
$q4.Cells.Item(22,2).Value = "009608"
$q4.Cells.Item(22,3).Value = "广发中证500指数增强A"
$q4.Cells.Item(22,4).Value = "1.17"
$q4.Cells.Item(22,5).Value = "93.86"
$q4.Cells.Item(22,6).Value = "1.37"
$q4.Cells.Item(22,7).Value = "0.0160"
$q4.Cells.Item(22,8).Value = 5
$q4.Cells.Item(23,1).Value = 21
$q4.Cells.Item(23,2).Value = "009609"
$q4.Cells.Item(23,3).Value = "广发中证500指数增强C"
$q4.Cells.Item(23,4).Value = "0.75"
$q4.Cells.Item(23,5).Value = "93.86"
$q4.Cells.Item(23,6).Value = "1.37"
$q4.Cells.Item(23,7).Value = "0.0103"
$q4.Cells.Item(23,8).Value = 5
$q4.Cells.Item(24,1).Value = 22
$q4.Cells.Item(24,2).Value = "159918"
$q4.Cells.Item(24,3).Value = "嘉实中创400ETF"
$q4.Cells.Item(24,4).Value = "0.59"
$q4.Cells.Item(24,5).Value = "98.55"
$q4.Cells.Item(24,6).Value = "0.63"
$q4.Cells.Item(24,7).Value = "0.0037"
$q4.Cells.Item(24,8).Value = 10
$q4.Cells.Item(25,1).Value = 23
$q4.Cells.Item(25,2).Value = "006346"
$q4.Cells.Item(25,3).Value = "安信量化优选股票A"
$q4.Cells.Item(25,4).Value = "0.27"
$q4.Cells.Item(25,5).Value = "90.65"
$q4.Cells.Item(25,6).Value = "1.13"
$q4.Cells.Item(25,7).Value = "0.0031"
$q4.Cells.Item(25,8).Value = 8
$q4.Cells.Item(26,1).Value = 24
$q4.Cells.Item(26,2).Value = "006347"
$q4.Cells.Item(26,3).Value = "安信量化优选股票C"
$q4.Cells.Item(26,4).Value = "0.14"
$q4.Cells.Item(26,5).Value = "90.65"
$q4.Cells.Item(26,6).Value = "1.13"
$q4.Cells.Item(26,7).Value = "0.0016"
$q4.Cells.Item(26,8).Value = 8
$q4.Cells.Item(27,1).Value = 25
$q4.Cells.Item(27,2).Value = "006783"
$q4.Cells.Item(27,3).Value = "红土创新中证500指数增强A"
$q4.Cells.Item(27,4).Value = "0.05"
$q4.Cells.Item(27,5).Value = "92.80"
$q4.Cells.Item(27,6).Value = "2.19"
$q4.Cells.Item(27,7).Value = "0.0011"
$q4.Cells.Item(27,8).Value = 7
$q4.Cells.Item(28,1).Value = 26
$q4.Cells.Item(28,2).Value = "016209"
$q4.Cells.Item(28,3).Value = "申万菱信中证军工指数C"
$q4.Cells.Item(28,4).Value = "0.04"
$q4.Cells.Item(28,5).Value = "93.84"
$q4.Cells.Item(28,6).Value = "2.29"
$q4.Cells.Item(28,7).Value = "0.0009"
$q4.Cells.Item(28,8).Value = 10
$q4.Cells.Item(29,1).Value = 27
$q4.Cells.Item(29,2).Value = "006784"
$q4.Cells.Item(29,3).Value = "红土创新中证500指数增强C"
$q4.Cells.Item(29,4).Value = "0.04"
$q4.Cells.Item(29,5).Value = "92.80"
$q4.Cells.Item(29,6).Value = "2.19"
$q4.Cells.Item(29,7).Value = "0.0009"
$q4.Cells.Item(29,8).Value = 7

# Apply the same "index column" styling (bold/centered/bordered) used
# elsewhere in the workbook to column A of the new sheet's data rows.
$zj.Range("A2").Copy() | Out-Null
$q4.Range("A2:A29").PasteSpecial(-4122) | Out-Null
for ($r = 2; $r -le 29; $r++) {
    $q4.Cells.Item($r, 1).Value = $r - 2
}

# ---- Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the
# top of the data (row 2), pushing the existing quarters down by one, and
# renumber the index column to stay sequential. ----
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()
$zj.Range("A3").Copy() | Out-Null
$zj.Range("A2").PasteSpecial(-4122) | Out-Null

$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q4"
$zj.Cells.Item(2,3).Value = 28
$zj.Cells.Item(2,4).Value = 6.96

$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(4,1).Value = 2
$zj.Cells.Item(5,1).Value = 3
$zj.Cells.Item(6,1).Value = 4
